$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update rotates the Fecha/Volumen/Precio data down one row
# (row 5's values move up to row 2, and rows 2-4 each shift down one row).

$ws.Range("D2").Value = 44277
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("P2").Value = 550

$ws.Range("D3").Value = 44284

$ws.Range("D4").Value = 44280
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 500

$ws.Range("D5").Value = 44291
$ws.Range("J5").Value = 30
